$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet after the last existing sheet (i.e. after "Pythoncode")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Sheet2"

# Header row: same formatting as sheet1's header (copy the format only),
# then fill in the same header text.
$ws1.Range("A1:B1").Copy()
$ws2.Range("A1:B1").PasteSpecial(-4122)
$ws2.Range("A1").Value = "pythonCode"
$ws2.Range("B1").Value = "Result"

# Data row (plain, unstyled cells)
$ws2.Range("A2").Value = "Python is crap"
$ws2.Range("B2").Value = "NameError: name 'Python' is not defined on line 1"

# Update selections: sheet1 moves its selection to A3:B3, and the new
# Sheet2 becomes the active sheet/tab with A2:B2 selected.
$ws1.Range("A3:B3").Select()
$ws2.Range("A2:B2").Select()
